# correção nos dados e inicio da analise PNAD 2009
#
# The two "section header" rows (row 5 "situação do domicílio" and row 8
# "grandes regiões e unidades da federação") were removed from the sheet;
# everything below each of them shifts up to close the gap. The row-2
# placeholder labels ("unnamed: 1_level_1" / "unnamed: 5_level_1") are
# corrected to "total" to match the neighbouring "total" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled header cells on row 2.
$ws.Range("B2").Value2 = "total"
$ws.Range("F2").Value2 = "total"

# Remove the blank "grandes regiões e unidades da federação" section-header
# row first (higher row number first so row numbers above it don't shift).
$ws.Rows(8).Delete()

# Remove the blank "situação do domicílio" section-header row.
$ws.Rows(5).Delete()
